$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style) from H1 (existing header cell) to the new header cells I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Set header values for the new columns
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Set data values for the new columns
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8

$ws.Range("I3").Value = 8
$ws.Range("J3").Value = 8
